$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6206553381867081
$ws.Range("C2").Value = 0.3267770552806262
$ws.Range("E2").Value = 0.2783447425719281
$ws.Range("F2").Value = 1.945874857898332
$ws.Range("G2").Value = 0.6882202465645761
$ws.Range("H2").Value = 0.8093703542944155
$ws.Range("J2").Value = 0.05405308978529888
$ws.Range("K2").Value = 0.285285936532631
$ws.Range("L2").Value = 0.4570528690655919
$ws.Range("O2").Value = 2.992967561356707
$ws.Range("B3").Value = 0.5744671232763494
$ws.Range("C3").Value = 0.3293776342582859
$ws.Range("E3").Value = 0.2770727828316808
$ws.Range("F3").Value = 1.948879681086069
$ws.Range("G3").Value = 0.6950705483605475
$ws.Range("H3").Value = 0.8163355163702661
$ws.Range("J3").Value = 0.05245750387063808
$ws.Range("K3").Value = 0.2496342426151443
$ws.Range("L3").Value = 0.4457837983816972
$ws.Range("O3").Value = 3.021829277669667
$ws.Range("B4").Value = 0.5462013731338971
$ws.Range("C4").Value = 0.3310632386890475
$ws.Range("E4").Value = 0.2764002169633955
$ws.Range("F4").Value = 1.951701272475148
$ws.Range("G4").Value = 0.6996966688780901
$ws.Range("H4").Value = 0.8209312980085173
$ws.Range("J4").Value = 0.05147243164009652
$ws.Range("K4").Value = 0.2276778722190755
$ws.Range("L4").Value = 0.4390296521563215
$ws.Range("O4").Value = 3.041100718047886
$ws.Range("B5").Value = 0.5347073514663236
$ws.Range("C5").Value = 0.3317725080468126
$ws.Range("E5").Value = 0.2761534821038509
$ws.Range("F5").Value = 1.953096923677791
$ws.Range("G5").Value = 0.7016874233062751
$ws.Range("H5").Value = 0.8228844266990123
$ws.Range("J5").Value = 0.05106967825046382
$ws.Range("K5").Value = 0.2187144606285187
$ws.Range("L5").Value = 0.4363190541772468
$ws.Range("O5").Value = 3.049343834289488
$ws.Range("B6").Value = 0.5328002831919036
$ws.Range("C6").Value = 0.3318916338550206
$ws.Range("E6").Value = 0.2761141655866055
$ws.Range("F6").Value = 1.95334352539524
$ws.Range("G6").Value = 0.702024361676358
$ws.Range("H6").Value = 0.8232135945154795
$ws.Range("J6").Value = 0.05100272180939314
$ws.Range("K6").Value = 0.2172251437216772
$ws.Range("L6").Value = 0.4358714918017199
$ws.Range("O6").Value = 3.050736144197501
$ws.Range("B7").Value = 0.5460462603866176
$ws.Range("C7").Value = 0.3310727135159013
$ws.Range("E7").Value = 0.2763967786043331
$ws.Range("F7").Value = 1.951719099026548
$ws.Range("G7").Value = 0.6997230894837827
$ws.Range("H7").Value = 0.8209573133336008
$ws.Range("J7").Value = 0.0514670053152102
$ws.Range("K7").Value = 0.2275570525450945
$ws.Range("L7").Value = 0.4389929266160522
$ws.Range("O7").Value = 3.041210309078636
$ws.Range("B8").Value = 0.6047106346291571
$ws.Range("C8").Value = 0.3276553163445488
$ws.Range("E8").Value = 0.2778837132711516
$ws.Range("F8").Value = 1.946708370792066
$ws.Range("G8").Value = 0.6904950331436552
$ws.Range("H8").Value = 0.811705739577441
$ws.Range("J8").Value = 0.05350406113381112
$ws.Range("K8").Value = 0.2730073005281497
$ws.Range("L8").Value = 0.4531331729395731
$ws.Range("O8").Value = 3.002597349552317
$ws.Range("B9").Value = 0.7204650701385447
$ws.Range("C9").Value = 0.321657167297392
$ws.Range("E9").Value = 0.2816568683676763
$ws.Range("F9").Value = 1.944622227773266
$ws.Range("G9").Value = 0.6757334140824369
$ws.Range("H9").Value = 0.7960930197198977
$ws.Range("J9").Value = 0.05745521526228003
$ws.Range("K9").Value = 0.3615883772946518
$ws.Range("L9").Value = 0.4821628156535098
$ws.Range("O9").Value = 2.939176723485659
$ws.Range("B10").Value = 0.8059105286681927
$ws.Range("C10").Value = 0.3176768995585313
$ws.Range("E10").Value = 0.284948178536645
$ws.Range("F10").Value = 1.947798491924985
$ws.Range("G10").Value = 0.6669238491234069
$ws.Range("H10").Value = 0.7861607397058705
$ws.Range("J10").Value = 0.06033072540976292
$ws.Range("K10").Value = 0.4263118002328383
$ws.Range("L10").Value = 0.5042737372267254
$ws.Range("O10").Value = 2.90007838213576
$ws.Range("B11").Value = 0.8448621244457968
$ws.Range("C11").Value = 0.315958328200372
$ws.Range("E11").Value = 0.2865575194864718
$ws.Range("F11").Value = 1.950264128018517
$ws.Range("G11").Value = 0.6633588554547671
$ws.Range("H11").Value = 0.7819755529435781
$ws.Range("J11").Value = 0.0616327472216156
$ws.Range("K11").Value = 0.4556739814019863
$ws.Range("L11").Value = 0.5145005725524356
$ws.Range("O11").Value = 2.883919213253009
$ws.Range("B12").Value = 0.8596231178937614
$ws.Range("C12").Value = 0.3153207567278216
$ws.Range("E12").Value = 0.2871829858600208
$ws.Range("F12").Value = 1.951344375344064
$ws.Range("G12").Value = 0.6620725733127983
$ws.Range("H12").Value = 0.7804385693877265
$ws.Range("J12").Value = 0.06212489657890075
$ws.Range("K12").Value = 0.4667805251117727
$ws.Range("L12").Value = 0.5183972113485709
$ws.Range("O12").Value = 2.878034101683099
$ws.Range("B13").Value = 0.8564436051426014
$ws.Range("C13").Value = 0.3154574818524587
$ws.Range("E13").Value = 0.2870475681570497
$ws.Range("F13").Value = 1.951105209724261
$ws.Range("G13").Value = 0.662346762852188
$ws.Range("H13").Value = 0.7807674584870625
$ws.Range("J13").Value = 0.06201894376909678
$ws.Range("K13").Value = 0.4643890886132169
$ws.Range("L13").Value = 0.517556939169225
$ws.Range("O13").Value = 2.879291156706572
$ws.Range("B14").Value = 0.8460763077160891
$ws.Range("C14").Value = 0.3159056101240321
$ws.Range("E14").Value = 0.2866086558796184
$ws.Range("F14").Value = 1.950350064491317
$ws.Range("G14").Value = 0.6632517551382691
$ws.Range("H14").Value = 0.7818481454920345
$ws.Range("J14").Value = 0.06167325479410835
$ws.Range("K14").Value = 0.4565879732246287
$ws.Range("L14").Value = 0.514820672742573
$ws.Range("O14").Value = 2.883430351151446
$ws.Range("B15").Value = 0.839727426777813
$ws.Range("C15").Value = 0.3161818217525418
$ws.Range("E15").Value = 0.2863418965360935
$ws.Range("F15").Value = 1.949906596800275
$ws.Range("G15").Value = 0.6638143871844377
$ws.Range("H15").Value = 0.7825163286255972
$ws.Range("J15").Value = 0.06146139258227379
$ws.Range("K15").Value = 0.4518079465953235
$ws.Range("L15").Value = 0.5131477426880622
$ws.Range("O15").Value = 2.885996205111653
$ws.Range("B16").Value = 0.8033665216762813
$ws.Range("C16").Value = 0.3177910627664637
$ws.Range("E16").Value = 0.2848452536261021
$ws.Range("F16").Value = 1.947657878197901
$ws.Range("G16").Value = 0.667165740467361
$ws.Range("H16").Value = 0.7864409502633123
$ws.Range("J16").Value = 0.0602455109873361
$ws.Range("K16").Value = 0.4243912301967327
$ws.Range("L16").Value = 0.5036087566336818
$ws.Range("O16").Value = 2.9011671648137
$ws.Range("B17").Value = 0.7810806338040663
$ws.Range("C17").Value = 0.3188018435612641
$ws.Range("E17").Value = 0.2839557702982844
$ws.Range("F17").Value = 1.946539644108725
$ws.Range("G17").Value = 0.6693350772265134
$ws.Range("H17").Value = 0.788933855765066
$ws.Range("J17").Value = 0.05949803563024858
$ws.Range("K17").Value = 0.4075507944153287
$ws.Range("L17").Value = 0.4977998559286902
$ws.Range("O17").Value = 2.910890788732189
$ws.Range("B18").Value = 0.7682701566120613
$ws.Range("C18").Value = 0.3193918857437197
$ws.Range("E18").Value = 0.2834547197081392
$ws.Range("F18").Value = 1.945992555127319
$ws.Range("G18").Value = 0.6706244704659241
$ws.Range("H18").Value = 0.7903990602107029
$ws.Range("J18").Value = 0.05906753854855396
$ws.Range("K18").Value = 0.397857052127307
$ws.Range("L18").Value = 0.4944746026899338
$ws.Range("O18").Value = 2.916636683656165
$ws.Range("B19").Value = 0.7639341129555248
$ws.Range("C19").Value = 0.3195931533792926
$ws.Range("E19").Value = 0.28328688805572
$ws.Range("F19").Value = 1.945823830263578
$ws.Range("G19").Value = 0.6710681878377116
$ws.Range("H19").Value = 0.790900538993462
$ws.Range("J19").Value = 0.05892168270420939
$ws.Range("K19").Value = 0.3945736375306694
$ws.Range("L19").Value = 0.4933514633831777
$ws.Range("O19").Value = 2.918608441813916
$ws.Range("B20").Value = 0.7834522065339513
$ws.Range("C20").Value = 0.3186933473220002
$ws.Range("E20").Value = 0.2840493652967879
$ws.Range("F20").Value = 1.946648738676643
$ws.Range("G20").Value = 0.6690998366959136
$ws.Range("H20").Value = 0.7886652372649721
$ws.Range("J20").Value = 0.05957766474059412
$ws.Range("K20").Value = 0.4093442756657737
$ws.Range("L20").Value = 0.4984165824842819
$ws.Range("O20").Value = 2.909839844104766
$ws.Range("B21").Value = 0.8491211466507593
$ws.Range("C21").Value = 0.3157736255609773
$ws.Range("E21").Value = 0.2867371402428631
$ws.Range("F21").Value = 1.950567892966788
$ws.Range("G21").Value = 0.6629842075185124
$ws.Range("H21").Value = 0.781529423025134
$ws.Range("J21").Value = 0.06177481660844819
$ws.Range("K21").Value = 0.458879688310958
$ws.Range("L21").Value = 0.515623732430015
$ws.Range("O21").Value = 2.882208217720887
$ws.Range("B22").Value = 0.8921025228650024
$ws.Range("C22").Value = 0.3139424294328155
$ws.Range("E22").Value = 0.2885872367708515
$ws.Range("F22").Value = 1.953983475071325
$ws.Range("G22").Value = 0.6593586353923797
$ws.Range("H22").Value = 0.7771446792531123
$ws.Range("J22").Value = 0.0632055356413872
$ws.Range("K22").Value = 0.4911820822560742
$ws.Range("L22").Value = 0.5270091224218305
$ws.Range("O22").Value = 2.865513404844876
$ws.Range("B23").Value = 0.8691571009333074
$ws.Range("C23").Value = 0.3149127351718199
$ws.Range("E23").Value = 0.2875912774597609
$ws.Range("F23").Value = 1.952082428249895
$ws.Range("G23").Value = 0.6612596705681426
$ws.Range("H23").Value = 0.7794593902899152
$ws.Range("J23").Value = 0.06244242274471645
$ws.Range("K23").Value = 0.4739484774572702
$ws.Range("L23").Value = 0.5209198474187815
$ws.Range("O23").Value = 2.874298905503835
$ws.Range("B24").Value = 0.7823800122758371
$ws.Range("C24").Value = 0.3187423706499359
$ws.Range("E24").Value = 0.284007018830458
$ws.Range("F24").Value = 1.946599118602052
$ws.Range("G24").Value = 0.6692060574556891
$ws.Range("H24").Value = 0.7887865800085549
$ws.Range("J24").Value = 0.05954166679628514
$ws.Range("K24").Value = 0.408533479970572
$ws.Range("L24").Value = 0.4981377156564264
$ws.Range("O24").Value = 2.910314491323646
$ws.Range("B25").Value = 0.6890775865460625
$ws.Range("C25").Value = 0.3232047564961267
$ws.Range("E25").Value = 0.2805447084191712
$ws.Range("F25").Value = 1.944358953052912
$ws.Range("G25").Value = 0.6793695273701985
$ws.Range("H25").Value = 0.8000462502585819
$ws.Range("J25").Value = 0.05639107192389758
$ws.Range("K25").Value = 0.3376858499336493
$ws.Range("L25").Value = 0.4741713248005226
$ws.Range("O25").Value = 2.955017082043781
